# Delete the "Rajabhakshi Nadaf" employee row (row 56) from the sheet,
# shifting all rows below it up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Delete()
